$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room for a second data row: insert a blank row at position 5
#    (old row 4 = data row stays at 4; old row 5 = Source row shifts to 6;
#     the merged range A5:H5 automatically becomes A6:H6)
# ------------------------------------------------------------------
$ws.Rows.Item(5).Insert()

# ------------------------------------------------------------------
# 2. Row 1 - replace the title text, merge across A1:I1, taller row,
#    center/wrap alignment (style copied from the old title cell which
#    already carries the bold Arial font we need).
# ------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A1:I1").PasteSpecial(-4122) | Out-Null
$ws.Range("A1:I1").Merge() | Out-Null
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in C. Batumi Municipality"
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# ------------------------------------------------------------------
# 3. Row 2 - subtitle stays put, just drop back to the default row height
# ------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14.5

# ------------------------------------------------------------------
# 4. Row 3, cell A3 - font changes from Arial to Sylfaen (cell stays blank)
# ------------------------------------------------------------------
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.Family = 1
$ws.Range("A3").Font.Name = "Sylfaen"

# ------------------------------------------------------------------
# 5. Row 4 - becomes the first data line ("family with disabilities Persons")
#    Style copied from old A4 (same font/fill) then the bottom border is
#    dropped, leaving just the top border.
# ------------------------------------------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Borders.Item(9).LineStyle = -4142
$ws.Range("A4").Value = "family with disabilities Persons "

$ws.Range("C4").Copy() | Out-Null
$ws.Range("B4:I4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4:I4").Borders.Item(8).LineStyle = -4142
$ws.Range("B4:I4").Borders.Item(9).LineStyle = -4142
$ws.Range("B4").Value = 1446
$ws.Range("C4").Value = 1458
$ws.Range("D4").Value = 1501
$ws.Range("E4").Value = 1682
$ws.Range("F4").Value = 1813
$ws.Range("G4").Value = 1970
$ws.Range("H4").Value = 2082
$ws.Range("I4").Value = 2279
$ws.Rows.Item(4).RowHeight = 24.75

# ------------------------------------------------------------------
# 6. Row 5 (the newly inserted row) - second data line
#    ("disabilities Persons"). Style copied from old A4 then the top
#    border is dropped, leaving just the bottom border.
# ------------------------------------------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Borders.Item(8).LineStyle = -4142
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").Value = "disabilities Persons "

$ws.Range("C4").Copy() | Out-Null
$ws.Range("B5:I5").PasteSpecial(-4122) | Out-Null
$ws.Range("B5:I5").Borders.Item(8).LineStyle = -4142
$ws.Range("B5:I5").Borders.Item(9).LineStyle = -4142
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("B5").Value = 1678
$ws.Range("C5").Value = 1692
$ws.Range("D5").Value = 1731
$ws.Range("E5").Value = 1932
$ws.Range("F5").Value = 2071
$ws.Range("G5").Value = 2249
$ws.Range("H5").Value = 2397
$ws.Range("I5").Value = 2644
$ws.Rows.Item(5).RowHeight = 21

# ------------------------------------------------------------------
# 7. Row 6 (old Source row) - drop its top border on the label cell;
#    the blank merged cells (B6:H6) keep their existing top border.
# ------------------------------------------------------------------
$ws.Range("A6").Borders.Item(8).LineStyle = -4142
$ws.Rows.Item(6).RowHeight = 27.75

# ------------------------------------------------------------------
# 8. Column A width & final selection (matches where the author left off)
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.95
$ws.Range("A1:I1").Select() | Out-Null
